# Update the "dSF" column (F) values that were repulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -8
    9  = -1
    10 = -3
    13 = -5
    14 = 2
    18 = 9
    21 = -2
    24 = 1
    25 = 4
    26 = 2
    28 = -5
    31 = 5
    33 = 1
    36 = 3
    39 = -3
    40 = -1
    47 = -2
    48 = -1
    50 = -3
    51 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
